# Continue the "landing Page" regression block (rows 13-16) for 17 more rows
# (rows 17-33), cycling through the same 4 feature names (with an extra
# repeat of the first one on row 17) and incrementing the Developer/tester
# numbers, exactly as the existing rows 13-16 already do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Feature name used in each row, 17 through 33 (row 17 repeats the first
# name before the regular 4-item cycle resumes on row 18).
$names = @(
    "SELF PACED ONLINE TRAINING",
    "SELF PACED ONLINE TRAINING",
    "IN DEPTH MATERIAL",
    "LIFETIME INSTRUCTOR SUPPORT",
    "RESUME PREPARATION",
    "SELF PACED ONLINE TRAINING",
    "IN DEPTH MATERIAL",
    "LIFETIME INSTRUCTOR SUPPORT",
    "RESUME PREPARATION",
    "SELF PACED ONLINE TRAINING",
    "IN DEPTH MATERIAL",
    "LIFETIME INSTRUCTOR SUPPORT",
    "RESUME PREPARATION",
    "SELF PACED ONLINE TRAINING",
    "IN DEPTH MATERIAL",
    "LIFETIME INSTRUCTOR SUPPORT",
    "RESUME PREPARATION"
)

for ($row = 17; $row -le 33; $row++) {
    $name = $names[$row - 17]
    $num = $row - 1

    $ws.Range("A$row").Value = "validate $name"
    $ws.Range("B$row").Value = "Developer$num"
    $ws.Range("C$row").Value = "landing Page"
    $ws.Range("D$row").Value = "Validate $name"
    $ws.Range("E$row").Value = "1/4"
    $ws.Range("F$row").Value = "Pre-conditions"
    $ws.Range("G$row").Value = "N/A"
    $ws.Range("H$row").Value = "Open Site`nClick in New window`nValidar text $name"
    $ws.Range("I$row").Value = "Text should be in site"
    $ws.Range("J$row").Value = "tester$num"
    $ws.Range("K$row").Value = "Text isnt in front of the site"
    $ws.Range("L$row").Value = "Fail"
    $ws.Range("M$row").Value = "Test automation failed"
    $ws.Range("N$row").Value = "Or the output"

    # Writing the multi-line H column text causes the engine to auto-size
    # the row (customHeight). Reset it via AutoFit so the row keeps the
    # sheet's default height, matching rows 13-16.
    $ws.Rows.Item($row).AutoFit() | Out-Null
}
